$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of path-planning data appended after the existing A1:C2338 block.
# Each entry is Row, A, B, C (all numeric).
$newData = @(
    @(2339, 9.73745786833307, 9.737043825499235, 16.19082042100107),
    @(2340, 9.031789155793497, 9.031384911842718, 16.12696773335329),
    @(2341, 8.326120443253924, 8.325725998186201, 16.0631150457055),
    @(2342, 7.62045173071435, 7.620067084529683, 15.99926235805772),
    @(2343, 6.914783018174777, 6.914408170873166, 15.93540967040993),
    @(2344, 6.209114305635204, 6.208749257216649, 15.87155698276215),
    @(2345, 5.503445593095631, 5.503090343560132, 15.80770429511436),
    @(2346, 4.797776880556057, 4.797431429903615, 15.74385160746658),
    @(2347, 4.092108168016484, 4.091772516247099, 15.67999891981879),
    @(2348, 3.386439455476911, 3.386113602590582, 15.61614623217101),
    @(2349, 2.680770742937337, 2.680454688934065, 15.55229354452322),
    @(2350, 1.975102030397764, 1.974795775277549, 15.48844085687544),
    @(2351, 1.269433317858191, 1.269136861621032, 15.42458816922765),
    @(2352, 0.5637646053186174, 0.5634779479645151, 15.36073548157987),
    @(2353, -0.1419041072209559, -0.1421809656920016, 15.29688279393208),
    @(2354, -0.8475728197605292, -0.8478398793485182, 15.2330301062843),
    @(2355, -1.553241532300103, -1.553498793005035, 15.16917741863651),
    @(2356, -2.258910244839676, -2.259157706661552, 15.10532473098873),
    @(2357, -2.964578957379249, -2.964816620318068, 15.04147204334094),
    @(2358, -3.670247669918822, -3.670475533974585, 14.97761935569316),
    @(2359, -4.375916382458396, -4.376134447631102, 14.91376666804537),
    @(2360, -5.081585094997969, -5.081793361287619, 14.84991398039759),
    @(2361, -5.787253807537542, -5.787452274944135, 14.7860612927498),
    @(2362, -6.492922520077116, -6.493111188600652, 14.72220860510201),
    @(2363, -7.198591232616689, -7.198770102257169, 14.65835591745423),
    @(2364, -7.904259945156262, -7.904429015913685, 14.59450322980644),
    @(2365, -8.609928657695836, -8.610087929570202, 14.53065054215866),
    @(2366, -9.315597370235409, -9.315746843226719, 14.46679785451087),
    @(2367, -10.02126608277498, -10.02140575688324, 14.40294516686309),
    @(2368, -10.72693479531456, -10.72706467053975, 14.3390924792153),
    @(2369, -11.43260350785413, -11.43272358419627, 14.27523979156751),
    @(2370, -12.1382722203937, -12.13838249785279, 14.21138710391973),
    @(2371, -12.84394093293328, -12.84404141150931, 14.14753441627194),
    @(2372, -13.54960964547285, -13.54970032516582, 14.08368172862416),
    @(2373, -14.25527835801242, -14.25535923882234, 14.01982904097637),
    @(2374, -14.960947070552, -14.96101815247886, 13.95597635332859),
    @(2375, -15.66661578309157, -15.66667706613537, 13.8921236656808),
    @(2376, -16.37228449563114, -16.37233597979189, 13.82827097803302),
    @(2377, -17.07795320817072, -17.07799489344841, 13.76441829038523),
    @(2378, -17.78362192071029, -17.78365380710493, 13.70056560273744),
    @(2379, -18.48929063324986, -18.48931272076144, 13.63671291508966),
    @(2380, -19.19495934578944, -19.19497163441796, 13.57286022744187),
    @(2381, -19.90062805832901, -19.90063054807447, 13.50900753979409)
)

foreach ($entry in $newData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
}

Write-Host "Added" $newData.Count "rows; last row now:" $ws.Cells.Item(2381, 1).Value() $ws.Cells.Item(2381, 2).Value() $ws.Cells.Item(2381, 3).Value()
